# Auto-generated script applying the numeric corrections described in the commit diff.
# Each sheet in this workbook contains generated market-price/profit data (no formulas);
# the edit simply rewrites specific H:N cells with refreshed values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 485.08334
$ws.Range("J38").Value = 1111
$ws.Range("L38").Value = 3333
$ws.Range("N38").Value = -4077
# Row 40
$ws.Range("H40").Value = 2173
$ws.Range("I40").Value = 2712.8333
$ws.Range("J40").Value = 1768.125
$ws.Range("K40").Value = 2712.8333
$ws.Range("L40").Value = 1768.125
$ws.Range("M40").Value = -2537.8333
$ws.Range("N40").Value = -2118.125
# Row 64
$ws.Range("H64").Value = 3475.5334
$ws.Range("I64").Value = 3289.8696
$ws.Range("K64").Value = 3289.8696
$ws.Range("M64").Value = -3041.8696
# Row 67
$ws.Range("H67").Value = 3475.5334
$ws.Range("I67").Value = 3289.8696
$ws.Range("K67").Value = 3289.8696
$ws.Range("M67").Value = -2431.8696
# Row 69
$ws.Range("H69").Value = 5066.1113
$ws.Range("I69").Value = 5760.8335
$ws.Range("J69").Value = 3676.6667
$ws.Range("K69").Value = 17282.5005
$ws.Range("L69").Value = 11030.0001
$ws.Range("M69").Value = -16408.5005
$ws.Range("N69").Value = -12778.0001
# Row 72
$ws.Range("H72").Value = 5066.1113
$ws.Range("I72").Value = 5760.8335
$ws.Range("J72").Value = 3676.6667
$ws.Range("K72").Value = 51847.5015
$ws.Range("L72").Value = 33090.0003
$ws.Range("M72").Value = -47479.5015
$ws.Range("N72").Value = -41826.0003
# Row 76
$ws.Range("H76").Value = 3630.6943
$ws.Range("I76").Value = 3503.1667
$ws.Range("K76").Value = 3503.1667
$ws.Range("M76").Value = -3188.1667
# Row 79
$ws.Range("H79").Value = 3630.6943
$ws.Range("I79").Value = 3503.1667
$ws.Range("K79").Value = 3503.1667
$ws.Range("M79").Value = -2411.1667
# Row 111
$ws.Range("H111").Value = 5815.5
$ws.Range("I111").Value = 6920.6665
$ws.Range("K111").Value = 20761.9995
$ws.Range("M111").Value = -17694.9995
# Row 129
$ws.Range("H129").Value = 1137.25
$ws.Range("J129").Value = 1316.3334
$ws.Range("L129").Value = 3949.0002
$ws.Range("N129").Value = -13949.0002
# Row 132
$ws.Range("H132").Value = 1492.0339
$ws.Range("I132").Value = 1249.3541
$ws.Range("J132").Value = 2551
$ws.Range("K132").Value = 3748.0623
$ws.Range("L132").Value = 7653
$ws.Range("M132").Value = -1218.0623
$ws.Range("N132").Value = -12713
# Row 137
$ws.Range("H137").Value = 498145.53
$ws.Range("I137").Value = 1348.6061
$ws.Range("J137").Value = 1210941.1
$ws.Range("K137").Value = 4045.8183
$ws.Range("L137").Value = 3632823.3
$ws.Range("M137").Value = -1495.8183
$ws.Range("N137").Value = -3637923.3

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4594.575
$ws.Range("I61").Value = 3087.7144
$ws.Range("K61").Value = 3087.7144
$ws.Range("M61").Value = -2875.7144
# Row 63
$ws.Range("H63").Value = 1172
$ws.Range("I63").Value = 1172
$ws.Range("K63").Value = 1172
$ws.Range("M63").Value = -486
# Row 66
$ws.Range("H66").Value = 1172
$ws.Range("I66").Value = 1172
$ws.Range("K66").Value = 5860
$ws.Range("M66").Value = -2428
# Row 80
$ws.Range("H80").Value = 250023060
$ws.Range("I80").Value = 25000
$ws.Range("J80").Value = 333355740
$ws.Range("K80").Value = 25000
$ws.Range("L80").Value = 333355740
$ws.Range("M80").Value = -24002
$ws.Range("N80").Value = -333357736
# Row 83
$ws.Range("H83").Value = 250023060
$ws.Range("I83").Value = 25000
$ws.Range("J83").Value = 333355740
$ws.Range("K83").Value = 75000
$ws.Range("L83").Value = 1000067220
$ws.Range("M83").Value = -70008
$ws.Range("N83").Value = -1000077204
# Row 122
$ws.Range("H122").Value = 7355681
$ws.Range("I122").Value = 6131
$ws.Range("J122").Value = 9617081
$ws.Range("K122").Value = 18393
$ws.Range("L122").Value = 28851243
$ws.Range("M122").Value = -15943
$ws.Range("N122").Value = -28856143
# Row 132
$ws.Range("H132").Value = 3756.228
$ws.Range("I132").Value = 1138.3414
$ws.Range("J132").Value = 10464.5625
$ws.Range("K132").Value = 3415.0242
$ws.Range("L132").Value = 31393.6875
$ws.Range("M132").Value = -885.0241999999998
$ws.Range("N132").Value = -36453.6875
# Row 136
$ws.Range("H136").Value = 4594.575
$ws.Range("I136").Value = 3087.7144
$ws.Range("K136").Value = 9263.143199999999
$ws.Range("M136").Value = -6713.143199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1151.4706
$ws.Range("I20").Value = 1165.9231
$ws.Range("J20").Value = 1104.5
$ws.Range("K20").Value = 1165.9231
$ws.Range("L20").Value = 1104.5
$ws.Range("M20").Value = -918.9231
$ws.Range("N20").Value = -1598.5
# Row 25
$ws.Range("H25").Value = 1328
$ws.Range("I25").Value = 1328
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1328
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1093
$ws.Range("N25").ClearContents()
# Row 57
$ws.Range("H57").Value = 60000
$ws.Range("J57").Value = 60000
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61440
# Row 94
$ws.Range("H94").Value = 1505.4445
$ws.Range("I94").Value = 1462.8462
$ws.Range("K94").Value = 1462.8462
$ws.Range("M94").Value = -1011.8462
# Row 134
$ws.Range("H134").Value = 4676.3667
$ws.Range("I134").Value = 4040.5
$ws.Range("J134").Value = 6425
$ws.Range("K134").Value = 12121.5
$ws.Range("L134").Value = 19275
$ws.Range("M134").Value = -9586.5
$ws.Range("N134").Value = -24345
# Row 136
$ws.Range("H136").Value = 60000
$ws.Range("J136").Value = 60000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2013.5306
$ws.Range("J31").Value = 6219
$ws.Range("L31").Value = 6219
$ws.Range("N31").Value = -6809
# Row 34
$ws.Range("H34").Value = 2013.5306
$ws.Range("J34").Value = 6219
$ws.Range("L34").Value = 6219
$ws.Range("N34").Value = -6623
# Row 44
$ws.Range("H44").Value = 23221
$ws.Range("J44").Value = 29799.5
$ws.Range("L44").Value = 29799.5
$ws.Range("N44").Value = -30683.5
# Row 58
$ws.Range("H58").Value = 2333804.8
$ws.Range("I58").Value = 4330732.5
$ws.Range("J58").Value = 4056
$ws.Range("K58").Value = 4330732.5
$ws.Range("L58").Value = 4056
$ws.Range("M58").Value = -4330529.5
$ws.Range("N58").Value = -4462
# Row 134
$ws.Range("H134").Value = 3751.225
$ws.Range("I134").Value = 2829.6428
$ws.Range("J134").Value = 4247.4614
$ws.Range("K134").Value = 8488.928400000001
$ws.Range("L134").Value = 12742.3842
$ws.Range("M134").Value = -5953.928400000001
$ws.Range("N134").Value = -17812.3842
# Row 136
$ws.Range("H136").Value = 2333804.8
$ws.Range("I136").Value = 4330732.5
$ws.Range("J136").Value = 4056
$ws.Range("K136").Value = 12992197.5
$ws.Range("L136").Value = 12168
$ws.Range("M136").Value = -12989647.5
$ws.Range("N136").Value = -17268
# Row 141
$ws.Range("H141").Value = 48509.25
$ws.Range("J141").Value = 48431.566
$ws.Range("L141").Value = 48431.566
$ws.Range("N141").Value = -58791.566

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 55715390
$ws.Range("I9").Value = 3000
$ws.Range("J9").Value = 66857868
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 200573604
$ws.Range("M9").Value = -8776
$ws.Range("N9").Value = -200574052
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 80
$ws.Range("H80").Value = 5821.8887
$ws.Range("I80").Value = 8914.857
$ws.Range("J80").Value = 3853.6365
$ws.Range("K80").Value = 8914.857
$ws.Range("L80").Value = 3853.6365
$ws.Range("M80").Value = -7916.857
$ws.Range("N80").Value = -5849.636500000001
# Row 83
$ws.Range("H83").Value = 5821.8887
$ws.Range("I83").Value = 8914.857
$ws.Range("J83").Value = 3853.6365
$ws.Range("K83").Value = 44574.285
$ws.Range("L83").Value = 19268.1825
$ws.Range("M83").Value = -39582.285
$ws.Range("N83").Value = -29252.1825
# Row 122
$ws.Range("H122").Value = 6041.6665
$ws.Range("I122").Value = 9314.286
$ws.Range("K122").Value = 27942.858
$ws.Range("M122").Value = -25492.858
# Row 123
$ws.Range("H123").Value = 27494.875
$ws.Range("J123").Value = 27494.875
$ws.Range("L123").Value = 27494.875
$ws.Range("N123").Value = -32394.875
# Row 132
$ws.Range("H132").Value = 1224.8055
$ws.Range("I132").Value = 1005.96
$ws.Range("J132").Value = 1722.1818
$ws.Range("K132").Value = 3017.88
$ws.Range("L132").Value = 5166.5454
$ws.Range("M132").Value = -487.8800000000001
$ws.Range("N132").Value = -10226.5454

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 56
$ws.Range("H56").Value = 13840.714
$ws.Range("I56").Value = 9633.333000000001
$ws.Range("J56").Value = 14988.182
$ws.Range("K56").Value = 9633.333000000001
$ws.Range("L56").Value = 14988.182
$ws.Range("M56").Value = -8942.333000000001
$ws.Range("N56").Value = -16370.182
# Row 62
$ws.Range("H62").Value = 36249
$ws.Range("J62").Value = 36249
$ws.Range("L62").Value = 36249
$ws.Range("N62").Value = -37497
# Row 65
$ws.Range("H65").Value = 36249
$ws.Range("J65").Value = 36249
$ws.Range("L65").Value = 108747
$ws.Range("N65").Value = -114987

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# Row 24
$ws.Range("H24").Value = 10666.667
$ws.Range("I24").Value = 8000
$ws.Range("J24").Value = 12000
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = -7770
$ws.Range("N24").Value = -12460
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 132
$ws.Range("H132").Value = 1211.3556
$ws.Range("I132").Value = 1250.0625
$ws.Range("K132").Value = 3750.1875
$ws.Range("M132").Value = -1220.1875
